$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.455362044514542
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 4.358119930609447

# Row 3
$ws.Range("B3").Value = 0.6606524410359556
$ws.Range("C3").Value = 0.306821227259698
$ws.Range("D3").Value = 22.3905356188092
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 23.85224582316562
